$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data (final state) for columns D, M, N, O, P, R, S across rows 2-15.
# Column order within each row array: D, M, N, O, P, R, S
$data = @{
    2  = @(44592, 30,  8000, 8000, 8000, "Provincia de Linares", 4000)
    3  = @(44214, 48,  6000, 6000, 6000, "Provincia de Linares", 3000)
    4  = @(44614, 45,  6000, 6000, 6000, "Provincia de Linares", 3000)
    5  = @(44211, 45,  6000, 6000, 6000, "Provincia de Curicó", 3000)
    6  = @(44589, 60,  6000, 6000, 6000, "Provincia de Curicó", 3000)
    7  = @(44588, 160, 6500, 7000, 6750, "Provincia de Curicó", 3375)
    8  = @(44628, 40,  6000, 6000, 6000, "Provincia de Linares", 3000)
    9  = @(44585, 160, 6500, 7000, 6750, "Provincia de Curicó", 3375)
    10 = @(44582, 150, 6000, 6500, 6233, "Provincia de Curicó", 3116)
    11 = @(44209, 58,  6000, 6000, 6000, "Provincia de Curicó", 3000)
    12 = @(44627, 45,  6000, 6000, 6000, "Provincia de Linares", 3000)
    13 = @(44586, 80,  7000, 7000, 7000, "Provincia de Curicó", 3500)
    14 = @(44587, 165, 6500, 7000, 6742, "Provincia de Linares", 3371)
    15 = @(44606, 45,  7000, 7000, 7000, "Provincia de Linares", 3500)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value  = $vals[0]   # D - Fecha
    $ws.Cells.Item($row, 13).Value = $vals[1]   # M - Volumen
    $ws.Cells.Item($row, 14).Value = $vals[2]   # N - Precio minimo
    $ws.Cells.Item($row, 15).Value = $vals[3]   # O - Precio maximo
    $ws.Cells.Item($row, 16).Value = $vals[4]   # P - Precio promedio ponderado
    $ws.Cells.Item($row, 18).Value = $vals[5]   # R - Origen
    $ws.Cells.Item($row, 19).Value = $vals[6]   # S - Precio $/Kg
}
